$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-Style($c, $style) {
    switch ($style) {
        1 { $c.HorizontalAlignment = -4108; $c.VerticalAlignment = -4108; $c.WrapText = $false; $c.Font.Bold = $false }
        2 { $c.HorizontalAlignment = -4108; $c.VerticalAlignment = -4108; $c.WrapText = $false; $c.Font.Bold = $true }
        3 { $c.HorizontalAlignment = -4108; $c.VerticalAlignment = -4108; $c.WrapText = $true;  $c.Font.Bold = $false }
        4 { $c.HorizontalAlignment = -4108; $c.VerticalAlignment = -4108; $c.WrapText = $true;  $c.Font.Bold = $true }
        5 { $c.HorizontalAlignment = -4131; $c.VerticalAlignment = -4108; $c.WrapText = $false; $c.Font.Bold = $false }
    }
}

function Set-Cell($addr, $style, $text) {
    $c = $ws.Range($addr)
    $c.Value2 = $text
    Set-Style $c $style
}

# --- Table 1: "Réalisation du système" ------------------------------------
Set-Cell "B2" 2 "Réalisation du système"

Set-Cell "C3" 2 "Tâche n°"
Set-Cell "D3" 4 "Nom de la tâche"
Set-Cell "E3" 2 "Temps nécessaire (projection)"
Set-Cell "F3" 2 "Réalisée"
Set-Cell "G3" 2 "Information potentielle"

Set-Cell "D4" 3 "Mettre en place son environnement de travail"
Set-Cell "E4" 1 "-"
Set-Cell "F4" 1 "Oui"

Set-Cell "D5" 3 "Choisir un language de programmation"
Set-Cell "E5" 1 "-"
Set-Cell "F5" 1 "Oui"
Set-Cell "G5" 1 "JavaScript"

Set-Cell "D6" 3 "Sélectionner un logiciel de développement adapté"
Set-Cell "E6" 1 "-"
Set-Cell "F6" 1 "Oui"
Set-Cell "G6" 1 "Visual Studio Code"

# Rows 7 & 8: tasks removed -> D cleared entirely (formula-driven C stays "-")
$ws.Range("D7").Clear()
Set-Cell "E7" 1 "-"

$ws.Range("D8").Clear()
Set-Cell "E8" 1 "-"

Set-Cell "D9" 3 "Réaliser un diagramme de séquence "
Set-Cell "E9" 1 "-"

Set-Cell "D10" 3 "Créer la base de donnée"
Set-Cell "E10" 1 "-"

Set-Cell "D11" 3 "Mettre en place l'interface de la page d'accueil"
Set-Cell "E11" 1 "-"

Set-Cell "D12" 3 "Coder le programme"
Set-Cell "E12" 1 "-"

Set-Cell "D13" 3 "Optimiser le code"
Set-Cell "E13" 1 "-"

# --- Table 2: "Test du système" --------------------------------------------
Set-Cell "B14" 2 "Test du système"

Set-Cell "D16" 3 "Préparation de la fiche recette"
Set-Cell "E16" 1 "-"

Set-Cell "D17" 3 "Vérification matériel"
Set-Cell "E17" 1 "-"

Set-Cell "D18" 3 "Mise en place du matériel (ordinateur, afficheur, etc…)"
Set-Cell "E18" 1 "-"

Set-Cell "D19" 3 "Test de l'application (fonctionnement bouton, bdd, envoi trame, etc…) + correction si nécessaire"
Set-Cell "E19" 1 "-"

Set-Cell "D20" 3 "Communication entre les appareils"
Set-Cell "E20" 1 "-"

# New recette / documentation tasks (rows 27-30)
Set-Cell "D27" 3 "établir un cahier de recette"
Set-Cell "D28" 3 "générer la documentation du code"
Set-Cell "D29" 3 "réaliser la documentation d'installation"
Set-Cell "D30" 3 "réaliser la documentation utilisateur"

Set-Cell "D32" 3 "Coder le programme"

Set-Cell "E33" 1 "implémenter le bouton envoyer indice"
Set-Cell "E34" 1 "implémenter le bouton modifier indice"
Set-Cell "E35" 1 "…"
Set-Cell "E36" 1 "gérer le test de connexion à l'afficheur"
Set-Cell "E37" 1 "afficher la liste des indices"
Set-Cell "E38" 1 "créer une première page"
Set-Cell "E39" 1 "installer WAMP"

Set-Cell "E40" 3 "Mettre en place l'interface de visualisation des indices"
Set-Cell "E41" 3 "Mettre en place l'interface d'aide"

Set-Cell "E42" 1 "implémenter la connexion à la BDD"
Set-Cell "E43" 1 "envoyer le texte ""Bonjour"" à l'afficheur"
Set-Cell "E44" 1 "paramétrer le convertisseur "
Set-Cell "E45" 1 "construire une trame à partir des paramètres (police, …)"

Set-Cell "E47" 3 "diagrammes de séquences"
Set-Cell "F48" 5 "diagramme de séquence de ""envoyer un indice"""
Set-Cell "F49" 5 "diagramme de séquence de ""ajouter un indice"""

Set-Cell "D50" 3 "valider le choix du javascript"
Set-Cell "E51" 1 "comment en javascript envoyer un message sur un port com"
Set-Cell "D52" 3 "valider le schéma réseau"

# --- View state (best effort) ----------------------------------------------
$ws.Range("E35").Select()
$excel.ActiveWindow.ScrollRow = 6
$excel.ActiveWindow.ScrollColumn = 1
